$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.161.98"
$ws.Range("E2").Value = "  +6.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.717.77"
$ws.Range("E3").Value = "  +3.88%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.20"
$ws.Range("E5").Value = "  +4.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3696"
$ws.Range("E7").Value = "  +1.93%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.29"
$ws.Range("E8").Value = "  +4.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3344"
$ws.Range("E9").Value = "  +2.50%  "
$ws.Range("E10").Value = "  +4.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07469"
$ws.Range("E11").Value = "  +5.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9999"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.335"
$ws.Range("E13").Value = "  +5.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.13"
$ws.Range("E14").Value = "  +2.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.943"
$ws.Range("E15").Value = "  +4.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.729.41"
$ws.Range("E16").Value = "  +3.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001079"
$ws.Range("E17").Value = "  +3.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06652"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "81.96"
$ws.Range("E19").Value = "  +3.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9990"
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.41"
$ws.Range("E21").Value = "  +4.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.090"
$ws.Range("E22").Value = "  +2.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.04"
$ws.Range("E23").Value = "  +3.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "26.098.75"
$ws.Range("E24").Value = "  +5.60%  "
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.453"
$ws.Range("E26").Value = "  +2.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.87"
$ws.Range("E27").Value = "  +2.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.374"
$ws.Range("E28").Value = "  +13.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.26"
$ws.Range("E29").Value = "  +3.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.913.94"
$ws.Range("E30").Value = "  +3.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "129.44"
$ws.Range("E31").Value = "  +2.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.107"
$ws.Range("E32").Value = "  +0.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.962"
$ws.Range("E33").Value = "  +1.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08563"
$ws.Range("E34").Value = "  +1.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.711"
$ws.Range("E35").Value = "  +3.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.92"
$ws.Range("E36").Value = "  +5.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.363"
$ws.Range("E37").Value = "  +2.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02332"
$ws.Range("E38").Value = "  +4.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06239"
$ws.Range("E39").Value = "  +3.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.590"
$ws.Range("E40").Value = "  +4.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2140"
$ws.Range("E41").Value = "  +3.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.235"
$ws.Range("E42").Value = "  -3.06%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6179"
$ws.Range("E43").Value = "  +4.25%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.44"
$ws.Range("E44").Value = "  +12.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9998"
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.836"
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5910"
$ws.Range("E47").Value = "  +4.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.76"
$ws.Range("E48").Value = "  +2.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.022"
$ws.Range("E49").Value = "  +3.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07259"
$ws.Range("E50").Value = "  +4.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.09"
$ws.Range("E51").Value = "  +3.34%  "
